$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '28.430.99'
$ws.Range("E2").Value = '  -0.31%  '

# Row 3
$ws.Range("D3").Value = '1.823.09'
$ws.Range("E3").Value = '  -0.52%  '

# Row 4
$ws.Range("E4").Value = '  +0.18%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '314.43'
$ws.Range("E5").Value = '  -0.87%  '

# Row 6
$ws.Range("E6").Value = '  +0.10%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5107'
$ws.Range("E7").Value = '  -3.71%  '

# Row 8
$ws.Range("E8").Value = '  -3.09%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07677'
$ws.Range("E9").Value = '  +1.20%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.93'
$ws.Range("E10").Value = '  -0.05%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.105'
$ws.Range("E11").Value = '  -0.52%  '

# Row 12
$ws.Range("E12").Value = '  +0.47%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.273'
$ws.Range("E13").Value = '  -0.93%  '

# Row 14
$ws.Range("E14").Value = '  +0.24%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.503'
$ws.Range("E15").Value = '  -0.95%  '

# Row 16
$ws.Range("D16").Value = '1.822.83'
$ws.Range("E16").Value = '  -0.55%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '92.87'
$ws.Range("E17").Value = '  +3.73%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001102'
$ws.Range("E18").Value = '  +2.77%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06639'
$ws.Range("E19").Value = '  +0.41%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.71'
$ws.Range("E20").Value = '  +0.69%  '

# Row 21
$ws.Range("E21").Value = '  +0.00%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.103'
$ws.Range("E22").Value = '  +0.56%  '

# Row 23
$ws.Range("D23").Value = '28.447.71'
$ws.Range("E23").Value = '  -0.35%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.16'
$ws.Range("E24").Value = '  -1.19%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.255'
$ws.Range("E25").Value = '  +5.96%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '21.10'
$ws.Range("E26").Value = '  +2.39%  '

# Row 27
$ws.Range("B27").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C27").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D27").Value = '2.032.71'
$ws.Range("E27").Value = '  -0.53%  '

# Row 28
$ws.Range("B28").Value = 'Monero'
$ws.Range("C28").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '155.80'
$ws.Range("E28").Value = '  -0.37%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.392'
$ws.Range("E29").Value = '  -3.74%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '124.05'
$ws.Range("E30").Value = '  +0.26%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.1099'
$ws.Range("E31").Value = '  +0.93%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.102'
$ws.Range("E32").Value = '  -2.29%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.649'
$ws.Range("E33").Value = '  -0.74%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.657'
$ws.Range("E34").Value = '  +0.08%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.07082'
$ws.Range("E35").Value = '  -1.49%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.2209'
$ws.Range("E36").Value = '  -2.64%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02324'
$ws.Range("E37").Value = '  -1.16%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.170'
$ws.Range("E38").Value = '  -1.53%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.768'
$ws.Range("E39").Value = '  +0.04%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6249'
$ws.Range("E40").Value = '  -0.57%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '11.18'
$ws.Range("E41").Value = '  -1.46%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.172'
$ws.Range("E42").Value = '  -1.26%  '

# Row 43
$ws.Range("E43").Value = '  +0.06%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.390'
$ws.Range("E44").Value = '  -1.12%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.31'
$ws.Range("E45").Value = '  -1.25%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.727'

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5861'
$ws.Range("E47").Value = '  -0.30%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '123.98'
$ws.Range("E48").Value = '  -1.66%  '

# Row 49
$ws.Range("E49").Value = '  -0.91%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.192'
$ws.Range("E50").Value = '  -0.19%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06899'
$ws.Range("E51").Value = '  +0.01%  '
